$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

# Row 2
Set-TextValue $ws.Range('D2') '36.345.95'
Set-TextValue $ws.Range('E2') '  +2.12%  '

# Row 3
Set-TextValue $ws.Range('D3') '2.006.48'
Set-TextValue $ws.Range('E3') '  +5.79%  '

# Row 4
Set-TextValue $ws.Range('E4') '  +0.02%  '

# Row 5
Set-TextValue $ws.Range('D5') '244.55'
Set-TextValue $ws.Range('E5') '  -1.00%  '

# Row 6
Set-TextValue $ws.Range('D6') '0.660'
Set-TextValue $ws.Range('E6') '  -4.47%  '

# Row 7
Set-TextValue $ws.Range('E7') '  +0.03%  '

# Row 8
Set-TextValue $ws.Range('D8') '44.42'
Set-TextValue $ws.Range('E8') '  +2.70%  '

# Row 9
Set-TextValue $ws.Range('D9') '61.88'
Set-TextValue $ws.Range('E9') '  +7.73%  '

# Row 10
Set-TextValue $ws.Range('D10') '0.363'
Set-TextValue $ws.Range('E10') '  +1.57%  '

# Row 11
Set-TextValue $ws.Range('D11') '0.0713'
Set-TextValue $ws.Range('E11') '  -5.66%  '

# Row 12
Set-TextValue $ws.Range('E12') '  -0.73%  '

# Row 13
Set-TextValue $ws.Range('D13') '14.36'
Set-TextValue $ws.Range('E13') '  -1.25%  '

# Row 14
Set-TextValue $ws.Range('E14') '  +5.59%  '

# Row 15
Set-TextValue $ws.Range('D15') '0.804'
Set-TextValue $ws.Range('E15') '  +0.13%  '

# Row 16
Set-TextValue $ws.Range('D16') '1.999.50'
Set-TextValue $ws.Range('E16') '  +5.90%  '

# Row 17
Set-TextValue $ws.Range('D17') '4.87'
Set-TextValue $ws.Range('E17') '  -3.12%  '

# Row 18
Set-TextValue $ws.Range('D18') '36.227.61'
Set-TextValue $ws.Range('E18') '  +1.77%  '

# Row 19
Set-TextValue $ws.Range('D19') '71.07'
Set-TextValue $ws.Range('E19') '  -3.69%  '

# Row 20
Set-TextValue $ws.Range('D20') '0.0₃0811'
Set-TextValue $ws.Range('E20') '  -2.64%  '

# Row 21
$ws.Range('B21').Value = 'BitcoinCash'
$ws.Range('C21').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
Set-TextValue $ws.Range('D21') '236.83'
Set-TextValue $ws.Range('E21') '  -3.76%  '

# Row 22
$ws.Range('B22').Value = 'Avalanche'
$ws.Range('C22').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
Set-TextValue $ws.Range('D22') '12.75'
Set-TextValue $ws.Range('E22') '  -1.96%  '

# Row 23
Set-TextValue $ws.Range('D23') '4.87'
Set-TextValue $ws.Range('E23') '  -6.35%  '

# Row 24
Set-TextValue $ws.Range('E24') '  +0.11%  '

# Row 25
Set-TextValue $ws.Range('D25') '2.43'
Set-TextValue $ws.Range('E25') '  -8.90%  '

# Row 26
Set-TextValue $ws.Range('D26') '165.33'
Set-TextValue $ws.Range('E26') '  -0.71%  '

# Row 27
Set-TextValue $ws.Range('D27') '8.60'
Set-TextValue $ws.Range('E27') '  -0.41%  '

# Row 28
Set-TextValue $ws.Range('D28') '19.43'
Set-TextValue $ws.Range('E28') '  +5.82%  '

# Row 29
Set-TextValue $ws.Range('D29') '1.93'
Set-TextValue $ws.Range('E29') '  -10.53%  '

# Row 30
Set-TextValue $ws.Range('D30') '0.121'
Set-TextValue $ws.Range('E30') '  -5.81%  '

# Row 31
Set-TextValue $ws.Range('D31') '21.72'
Set-TextValue $ws.Range('E31') '  +56.39%  '

# Row 32
Set-TextValue $ws.Range('D32') '4.33'
Set-TextValue $ws.Range('E32') '  -1.14%  '

# Row 33
Set-TextValue $ws.Range('D33') '0.0581'
Set-TextValue $ws.Range('E33') '  -3.47%  '

# Row 34
Set-TextValue $ws.Range('E34') '  +0.03%  '

# Row 35
$ws.Range('B35').Value = 'Kaspa'
$ws.Range('C35').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextValue $ws.Range('D35') '0.0870'
Set-TextValue $ws.Range('E35') '  +19.10%  '

# Row 36
$ws.Range('B36').Value = 'WEMIXToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
Set-TextValue $ws.Range('D36') '1.86'
Set-TextValue $ws.Range('E36') '  +0.73%  '

# Row 37
Set-TextValue $ws.Range('D37') '3.96'
Set-TextValue $ws.Range('E37') '  -6.99%  '

# Row 38
Set-TextValue $ws.Range('D38') '2.10'
Set-TextValue $ws.Range('E38') '  +7.07%  '

# Row 39
Set-TextValue $ws.Range('D39') '0.850'
Set-TextValue $ws.Range('E39') '  -0.59%  '

# Row 40
Set-TextValue $ws.Range('E40') '  -10.75%  '

# Row 41
Set-TextValue $ws.Range('D41') '0.0214'
Set-TextValue $ws.Range('E41') '  -5.42%  '

# Row 42
$ws.Range('B42').Value = 'Aave'
$ws.Range('C42').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue $ws.Range('D42') '95.22'
Set-TextValue $ws.Range('E42') '  -4.19%  '

# Row 43
$ws.Range('B43').Value = 'ARBITRUM'
$ws.Range('C43').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextValue $ws.Range('D43') '1.11'
Set-TextValue $ws.Range('E43') '  +1.88%  '

# Row 44
Set-TextValue $ws.Range('D44') '2.76'
Set-TextValue $ws.Range('E44') '  +15.47%  '

# Row 45
Set-TextValue $ws.Range('D45') '15.96'
Set-TextValue $ws.Range('E45') '  -6.85%  '

# Row 46
Set-TextValue $ws.Range('D46') '1.306.62'
Set-TextValue $ws.Range('E46') '  -1.21%  '

# Row 47
Set-TextValue $ws.Range('D47') '0.0811'
Set-TextValue $ws.Range('E47') '  -0.02%  '

# Row 48
Set-TextValue $ws.Range('D48') '2.77'
Set-TextValue $ws.Range('E48') '  +0.75%  '

# Row 49
$ws.Range('B49').Value = 'RenderToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue $ws.Range('D49') '2.19'
Set-TextValue $ws.Range('E49') '  -7.28%  '

# Row 50
$ws.Range('B50').Value = 'RocketPoolETH'
$ws.Range('C50').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
Set-TextValue $ws.Range('D50') '2.182.78'
Set-TextValue $ws.Range('E50') '  +5.16%  '

# Row 51
Set-TextValue $ws.Range('D51') '3.82'
Set-TextValue $ws.Range('E51') '  +14.42%  '
